$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs -> ECs
$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Csf2"
$ws.Cells.Item(2, 3).Value = "Csf2ra"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.2300786666666667
$ws.Cells.Item(2, 8).Value = 0.690236
$ws.Cells.Item(2, 9).Value = 0.2137022699341201
$ws.Cells.Item(2, 10).Value = 0.2304051507653011
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.5
$ws.Cells.Item(2, 13).Value = 0.2111075
$ws.Cells.Item(2, 14).Value = 0.422215
$ws.Cells.Item(2, 15).Value = 0.0009789996530386493
$ws.Cells.Item(2, 16).Value = 0.0006535518839479957
$ws.Cells.Item(2, 17).Value = 0.04857133212333333
$ws.Cells.Item(2, 18).Value = 0.29142799274
$ws.Cells.Item(2, 19).Value = 0.0002092144481190754
$ws.Cells.Item(2, 20).Value = 0.0001505817203539845

# Row 3: FAPs -> FAPs
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Csf2"
$ws.Cells.Item(3, 3).Value = "Csf2ra"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.2300786666666667
$ws.Cells.Item(3, 8).Value = 0.690236
$ws.Cells.Item(3, 9).Value = 0.2137022699341201
$ws.Cells.Item(3, 10).Value = 0.2304051507653011
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.5485233333333334
$ws.Cells.Item(3, 14).Value = 1.64557
$ws.Cells.Item(3, 15).Value = 0.00254374739418039
$ws.Cells.Item(3, 16).Value = 0.002547198402871294
$ws.Cells.Item(3, 17).Value = 0.1262035171688889
$ws.Cells.Item(3, 18).Value = 1.13583165452
$ws.Cells.Item(3, 19).Value = 0.0005436045922753524
$ws.Cells.Item(3, 20).Value = 0.0005868876320426946

# Row 4: FAPs -> Inflammatory-Mac
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Csf2"
$ws.Cells.Item(4, 3).Value = "Csf2ra"
$ws.Cells.Item(4, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.2300786666666667
$ws.Cells.Item(4, 8).Value = 0.690236
$ws.Cells.Item(4, 9).Value = 0.2137022699341201
$ws.Cells.Item(4, 10).Value = 0.2304051507653011
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 76.60934666666667
$ws.Cells.Item(4, 14).Value = 229.82804
$ws.Cells.Item(4, 15).Value = 0.3552717160981219
$ws.Cells.Item(4, 16).Value = 0.3557537001908395
$ws.Cells.Item(4, 17).Value = 17.62617633527111
$ws.Cells.Item(4, 18).Value = 158.63558701744
$ws.Cells.Item(4, 19).Value = 0.07592237217355895
$ws.Cells.Item(4, 20).Value = 0.08196748492778408

# Row 5: FAPs -> MuSCs
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Csf2"
$ws.Cells.Item(5, 3).Value = "Csf2ra"
$ws.Cells.Item(5, 4).Value = "MuSCs"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.2300786666666667
$ws.Cells.Item(5, 8).Value = 0.690236
$ws.Cells.Item(5, 9).Value = 0.2137022699341201
$ws.Cells.Item(5, 10).Value = 0.2304051507653011
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.6653395
$ws.Cells.Item(5, 14).Value = 1.330679
$ws.Cells.Item(5, 15).Value = 0.003085476070972885
$ws.Cells.Item(5, 16).Value = 0.00205977468204596
$ws.Cells.Item(5, 17).Value = 0.1530804250406667
$ws.Cells.Item(5, 18).Value = 0.9184825502439999
$ws.Cells.Item(5, 19).Value = 0.0006593732401943159
$ws.Cells.Item(5, 20).Value = 0.0004745826961593494

# Row 6: FAPs -> Neutrophils
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Csf2"
$ws.Cells.Item(6, 3).Value = "Csf2ra"
$ws.Cells.Item(6, 4).Value = "Neutrophils"
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 0.6666666666666666
$ws.Cells.Item(6, 7).Value = 0.2300786666666667
$ws.Cells.Item(6, 8).Value = 0.690236
$ws.Cells.Item(6, 9).Value = 0.2137022699341201
$ws.Cells.Item(6, 10).Value = 0.2304051507653011
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 98.13070699999999
$ws.Cells.Item(6, 14).Value = 294.392121
$ws.Cells.Item(6, 15).Value = 0.4550758646918624
$ws.Cells.Item(6, 16).Value = 0.4556932494084679
$ws.Cells.Item(6, 17).Value = 22.57778222561733
$ws.Cells.Item(6, 18).Value = 203.200040030556
$ws.Cells.Item(6, 19).Value = 0.09725074527688352
$ws.Cells.Item(6, 20).Value = 0.104994071832688

# Row 7: FAPs -> Resolving-Mac
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Csf2"
$ws.Cells.Item(7, 3).Value = "Csf2ra"
$ws.Cells.Item(7, 4).Value = "Resolving-Mac"
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 0.6666666666666666
$ws.Cells.Item(7, 7).Value = 0.2300786666666667
$ws.Cells.Item(7, 8).Value = 0.690236
$ws.Cells.Item(7, 9).Value = 0.2137022699341201
$ws.Cells.Item(7, 10).Value = 0.2304051507653011
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 39.47090533333333
$ws.Cells.Item(7, 14).Value = 118.412716
$ws.Cells.Item(7, 15).Value = 0.1830441960918238
$ws.Cells.Item(7, 16).Value = 0.1832925254318273
$ws.Cells.Item(7, 17).Value = 9.081413271219555
$ws.Cells.Item(7, 18).Value = 81.732719440976
$ws.Cells.Item(7, 19).Value = 0.03911696020308896
$ws.Cells.Item(7, 20).Value = 0.04223154195627295

# Row 8: MuSCs -> ECs
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 2).Value = "Csf2"
$ws.Cells.Item(8, 3).Value = "Csf2ra"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 0.5
$ws.Cells.Item(8, 7).Value = 0.2341465
$ws.Cells.Item(8, 8).Value = 0.468293
$ws.Cells.Item(8, 9).Value = 0.2174805655477089
$ws.Cells.Item(8, 10).Value = 0.1563191709318771
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.5
$ws.Cells.Item(8, 13).Value = 0.2111075
$ws.Cells.Item(8, 14).Value = 0.422215
$ws.Cells.Item(8, 15).Value = 0.0009789996530386493
$ws.Cells.Item(8, 16).Value = 0.0006535518839479957
$ws.Cells.Item(8, 17).Value = 0.04943008224875
$ws.Cells.Item(8, 18).Value = 0.197720328995
$ws.Cells.Item(8, 19).Value = 0.0002129133982138562
$ws.Cells.Item(8, 20).Value = 0.0001021626886597171

# Row 9: MuSCs -> FAPs
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 2).Value = "Csf2"
$ws.Cells.Item(9, 3).Value = "Csf2ra"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 0.5
$ws.Cells.Item(9, 7).Value = 0.2341465
$ws.Cells.Item(9, 8).Value = 0.468293
$ws.Cells.Item(9, 9).Value = 0.2174805655477089
$ws.Cells.Item(9, 10).Value = 0.1563191709318771
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.5485233333333334
$ws.Cells.Item(9, 14).Value = 1.64557
$ws.Cells.Item(9, 15).Value = 0.00254374739418039
$ws.Cells.Item(9, 16).Value = 0.002547198402871294
$ws.Cells.Item(9, 17).Value = 0.1284348186683333
$ws.Cells.Item(9, 18).Value = 0.77060891201
$ws.Cells.Item(9, 19).Value = 0.000553215621896862
$ws.Cells.Item(9, 20).Value = 0.0003981759425358422

# Row 10: MuSCs -> Inflammatory-Mac
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Csf2"
$ws.Cells.Item(10, 3).Value = "Csf2ra"
$ws.Cells.Item(10, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 0.5
$ws.Cells.Item(10, 7).Value = 0.2341465
$ws.Cells.Item(10, 8).Value = 0.468293
$ws.Cells.Item(10, 9).Value = 0.2174805655477089
$ws.Cells.Item(10, 10).Value = 0.1563191709318771
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 76.60934666666667
$ws.Cells.Item(10, 14).Value = 229.82804
$ws.Cells.Item(10, 15).Value = 0.3552717160981219
$ws.Cells.Item(10, 16).Value = 0.3557537001908395
$ws.Cells.Item(10, 17).Value = 17.93781038928667
$ws.Cells.Item(10, 18).Value = 107.62686233572
$ws.Cells.Item(10, 19).Value = 0.07726469374012462
$ws.Cells.Item(10, 20).Value = 0.05561112346977961

# Row 11: MuSCs -> MuSCs
$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Csf2"
$ws.Cells.Item(11, 3).Value = "Csf2ra"
$ws.Cells.Item(11, 4).Value = "MuSCs"
$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(11, 6).Value = 0.5
$ws.Cells.Item(11, 7).Value = 0.2341465
$ws.Cells.Item(11, 8).Value = 0.468293
$ws.Cells.Item(11, 9).Value = 0.2174805655477089
$ws.Cells.Item(11, 10).Value = 0.1563191709318771
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 0.6653395
$ws.Cells.Item(11, 14).Value = 1.330679
$ws.Cells.Item(11, 15).Value = 0.003085476070972885
$ws.Cells.Item(11, 16).Value = 0.00205977468204596
$ws.Cells.Item(11, 17).Value = 0.15578691523675
$ws.Cells.Item(11, 18).Value = 0.623147660947
$ws.Cells.Item(11, 19).Value = 0.0006710310808991059
$ws.Cells.Item(11, 20).Value = 0.0003219822706038953

# Row 12: MuSCs -> Neutrophils
$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "Csf2"
$ws.Cells.Item(12, 3).Value = "Csf2ra"
$ws.Cells.Item(12, 4).Value = "Neutrophils"
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = 0.5
$ws.Cells.Item(12, 7).Value = 0.2341465
$ws.Cells.Item(12, 8).Value = 0.468293
$ws.Cells.Item(12, 9).Value = 0.2174805655477089
$ws.Cells.Item(12, 10).Value = 0.1563191709318771
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 98.13070699999999
$ws.Cells.Item(12, 14).Value = 294.392121
$ws.Cells.Item(12, 15).Value = 0.4550758646918624
$ws.Cells.Item(12, 16).Value = 0.4556932494084679
$ws.Cells.Item(12, 17).Value = 22.9769615865755
$ws.Cells.Item(12, 18).Value = 137.861769519453
$ws.Cells.Item(12, 19).Value = 0.09897015642029888
$ws.Cells.Item(12, 20).Value = 0.07123359094678482

# Row 13: MuSCs -> Resolving-Mac
$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "Csf2"
$ws.Cells.Item(13, 3).Value = "Csf2ra"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 6).Value = 0.5
$ws.Cells.Item(13, 7).Value = 0.2341465
$ws.Cells.Item(13, 8).Value = 0.468293
$ws.Cells.Item(13, 9).Value = 0.2174805655477089
$ws.Cells.Item(13, 10).Value = 0.1563191709318771
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 39.47090533333333
$ws.Cells.Item(13, 14).Value = 118.412716
$ws.Cells.Item(13, 15).Value = 0.1830441960918238
$ws.Cells.Item(13, 16).Value = 0.1832925254318273
$ws.Cells.Item(13, 17).Value = 9.241974335631333
$ws.Cells.Item(13, 18).Value = 55.451846013788
$ws.Cells.Item(13, 19).Value = 0.03980855528627558
$ws.Cells.Item(13, 20).Value = 0.02865213561351325

# Row 14: Neutrophils -> ECs
$ws.Cells.Item(14, 1).Value = "Neutrophils"
$ws.Cells.Item(14, 2).Value = "Csf2"
$ws.Cells.Item(14, 3).Value = "Csf2ra"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 0.6124066666666667
$ws.Cells.Item(14, 8).Value = 1.83722
$ws.Cells.Item(14, 9).Value = 0.5688171645181709
$ws.Cells.Item(14, 10).Value = 0.6132756783028217
$ws.Cells.Item(14, 11).Value = 1
$ws.Cells.Item(14, 12).Value = 0.5
$ws.Cells.Item(14, 13).Value = 0.2111075
$ws.Cells.Item(14, 14).Value = 0.422215
$ws.Cells.Item(14, 15).Value = 0.0009789996530386493
$ws.Cells.Item(14, 16).Value = 0.0006535518839479957
$ws.Cells.Item(14, 17).Value = 0.1292836403833333
$ws.Cells.Item(14, 18).Value = 0.7757018422999999
$ws.Cells.Item(14, 19).Value = 0.0005568718067057176
$ws.Cells.Item(14, 20).Value = 0.0004008074749342941

# Row 15: Neutrophils -> FAPs
$ws.Cells.Item(15, 1).Value = "Neutrophils"
$ws.Cells.Item(15, 2).Value = "Csf2"
$ws.Cells.Item(15, 3).Value = "Csf2ra"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 0.6124066666666667
$ws.Cells.Item(15, 8).Value = 1.83722
$ws.Cells.Item(15, 9).Value = 0.5688171645181709
$ws.Cells.Item(15, 10).Value = 0.6132756783028217
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 0.5485233333333334
$ws.Cells.Item(15, 14).Value = 1.64557
$ws.Cells.Item(15, 15).Value = 0.00254374739418039
$ws.Cells.Item(15, 16).Value = 0.002547198402871294
$ws.Cells.Item(15, 17).Value = 0.3359193461555556
$ws.Cells.Item(15, 18).Value = 3.0232741154
$ws.Cells.Item(15, 19).Value = 0.001446927180008175
$ws.Cells.Item(15, 20).Value = 0.001562134828292757

# Row 16: Neutrophils -> Inflammatory-Mac
$ws.Cells.Item(16, 1).Value = "Neutrophils"
$ws.Cells.Item(16, 2).Value = "Csf2"
$ws.Cells.Item(16, 3).Value = "Csf2ra"
$ws.Cells.Item(16, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 0.6124066666666667
$ws.Cells.Item(16, 8).Value = 1.83722
$ws.Cells.Item(16, 9).Value = 0.5688171645181709
$ws.Cells.Item(16, 10).Value = 0.6132756783028217
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 76.60934666666667
$ws.Cells.Item(16, 14).Value = 229.82804
$ws.Cells.Item(16, 15).Value = 0.3552717160981219
$ws.Cells.Item(16, 16).Value = 0.3557537001908395
$ws.Cells.Item(16, 17).Value = 46.91607462764444
$ws.Cells.Item(16, 18).Value = 422.2446716488
$ws.Cells.Item(16, 19).Value = 0.2020846501844383
$ws.Cells.Item(16, 20).Value = 0.2181750917932757

# Row 17: Neutrophils -> MuSCs
$ws.Cells.Item(17, 1).Value = "Neutrophils"
$ws.Cells.Item(17, 2).Value = "Csf2"
$ws.Cells.Item(17, 3).Value = "Csf2ra"
$ws.Cells.Item(17, 4).Value = "MuSCs"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 0.6124066666666667
$ws.Cells.Item(17, 8).Value = 1.83722
$ws.Cells.Item(17, 9).Value = 0.5688171645181709
$ws.Cells.Item(17, 10).Value = 0.6132756783028217
$ws.Cells.Item(17, 11).Value = 2
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 0.6653395
$ws.Cells.Item(17, 14).Value = 1.330679
$ws.Cells.Item(17, 15).Value = 0.003085476070972885
$ws.Cells.Item(17, 16).Value = 0.00205977468204596
$ws.Cells.Item(17, 17).Value = 0.4074583453966666
$ws.Cells.Item(17, 18).Value = 2.44475007238
$ws.Cells.Item(17, 19).Value = 0.001755071749879463
$ws.Cells.Item(17, 20).Value = 0.001263209715282715

# Row 18: Neutrophils -> Neutrophils
$ws.Cells.Item(18, 1).Value = "Neutrophils"
$ws.Cells.Item(18, 2).Value = "Csf2"
$ws.Cells.Item(18, 3).Value = "Csf2ra"
$ws.Cells.Item(18, 4).Value = "Neutrophils"
$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 0.6124066666666667
$ws.Cells.Item(18, 8).Value = 1.83722
$ws.Cells.Item(18, 9).Value = 0.5688171645181709
$ws.Cells.Item(18, 10).Value = 0.6132756783028217
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 12).Value = 1
$ws.Cells.Item(18, 13).Value = 98.13070699999999
$ws.Cells.Item(18, 14).Value = 294.392121
$ws.Cells.Item(18, 15).Value = 0.4550758646918624
$ws.Cells.Item(18, 16).Value = 0.4556932494084679
$ws.Cells.Item(18, 17).Value = 60.09589917151332
$ws.Cells.Item(18, 18).Value = 540.8630925436199
$ws.Cells.Item(18, 19).Value = 0.25885496299468
$ws.Cells.Item(18, 20).Value = 0.2794655866289951

# Row 19: Neutrophils -> Resolving-Mac
$ws.Cells.Item(19, 1).Value = "Neutrophils"
$ws.Cells.Item(19, 2).Value = "Csf2"
$ws.Cells.Item(19, 3).Value = "Csf2ra"
$ws.Cells.Item(19, 4).Value = "Resolving-Mac"
$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 0.6124066666666667
$ws.Cells.Item(19, 8).Value = 1.83722
$ws.Cells.Item(19, 9).Value = 0.5688171645181709
$ws.Cells.Item(19, 10).Value = 0.6132756783028217
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 13).Value = 39.47090533333333
$ws.Cells.Item(19, 14).Value = 118.412716
$ws.Cells.Item(19, 15).Value = 0.1830441960918238
$ws.Cells.Item(19, 16).Value = 0.1832925254318273
$ws.Cells.Item(19, 17).Value = 24.17224556550222
$ws.Cells.Item(19, 18).Value = 217.55021008952
$ws.Cells.Item(19, 19).Value = 0.1041186806024593
$ws.Cells.Item(19, 20).Value = 0.1124088478620411
